# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 20:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 781368
$ws.Range("C4").Value = 16732
$ws.Range("E4").Value = 668023
$ws.Range("G4").Value = 1000
$ws.Range("H4").Value = 41575

# Row 25 - Israel
$ws.Range("B25").Value = 13713
$ws.Range("C25").Value = 222
$ws.Range("D25").Value = 4049
$ws.Range("E25").Value = 9487
$ws.Range("F25").Value = 149
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 177

# Row 32 - Rumania
$ws.Range("E32").Value = 6441
$ws.Range("G32").Value = 27
$ws.Range("H32").Value = 478

# Row 50 - Finlandia
$ws.Range("D50").Value = 2000
$ws.Range("E50").Value = 1770

# Row 64 - Barein
$ws.Range("B64").Value = 1907
$ws.Range("C64").Value = 26
$ws.Range("E64").Value = 1131

# Row 68 - Uzbekistan
$ws.Range("B68").Value = 1627
$ws.Range("C68").Value = 62
$ws.Range("E68").Value = 1361
